$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "43.500.10"
$ws.Range("E2").Value = "  -6.32%  "
Set-TextValue "D3" "2.535.39"
$ws.Range("E3").Value = "  -2.97%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "299.25"
$ws.Range("E5").Value = "  -2.96%  "
Set-TextValue "D6" "95.01"
$ws.Range("E6").Value = "  -5.37%  "
Set-TextValue "D7" "0.577"
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("E8").Value = "  +0.16%  "
Set-TextValue "D9" "0.556"
$ws.Range("E9").Value = "  -4.15%  "
Set-TextValue "D10" "36.70"
$ws.Range("E10").Value = "  -7.31%  "
Set-TextValue "D11" "0.0809"
$ws.Range("E11").Value = "  -4.40%  "
Set-TextValue "D12" "7.74"
$ws.Range("E12").Value = "  -5.30%  "
Set-TextValue "D13" "0.107"
$ws.Range("E13").Value = "  +1.12%  "
Set-TextValue "D14" "2.932.37"
$ws.Range("E14").Value = "  -2.57%  "
Set-TextValue "D15" "2.539.49"
$ws.Range("E15").Value = "  -2.67%  "
Set-TextValue "D16" "0.881"
$ws.Range("E16").Value = "  -4.37%  "
Set-TextValue "D17" "14.20"
$ws.Range("E17").Value = "  -4.93%  "
Set-TextValue "D18" "43.596.31"
$ws.Range("E18").Value = "  -6.41%  "
Set-TextValue "B19" "ShibaInu"
Set-TextValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.0₃0973"
$ws.Range("E19").Value = "  -3.99%  "
Set-TextValue "B20" "Uniswap"
Set-TextValue "C20" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.62"
$ws.Range("E20").Value = "  -1.78%  "
Set-TextValue "D21" "12.43"
$ws.Range("E21").Value = "  -4.26%  "
Set-TextValue "D22" "73.07"
$ws.Range("E22").Value = "  +2.14%  "
Set-TextValue "D23" "263.07"
$ws.Range("E23").Value = "  -3.99%  "
Set-TextValue "D24" "2.92"
$ws.Range("E24").Value = "  -4.18%  "
Set-TextValue "D25" "2.17"
$ws.Range("E25").Value = "  +0.47%  "
Set-TextValue "D26" "29.04"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.08%  "
Set-TextValue "D28" "10.13"
$ws.Range("E28").Value = "  -4.58%  "
Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.23"
$ws.Range("E29").Value = "  +0.37%  "
Set-TextValue "B30" "InjectiveProtocol"
Set-TextValue "C30" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "37.40"
$ws.Range("E30").Value = "  -4.33%  "
Set-TextValue "D31" "6.11"
$ws.Range("E31").Value = "  -4.06%  "
Set-TextValue "D32" "3.52"
$ws.Range("E32").Value = "  -3.73%  "
Set-TextValue "D33" "151.31"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  -3.18%  "
Set-TextValue "D36" "0.0804"
$ws.Range("E36").Value = "  -4.22%  "
Set-TextValue "D37" "0.116"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("E38").Value = "  -3.09%  "
Set-TextValue "D39" "23.68"
$ws.Range("E39").Value = "  +1.73%  "
Set-TextValue "D40" "16.36"
$ws.Range("E40").Value = "  +2.78%  "
Set-TextValue "D41" "3.53"
$ws.Range("E41").Value = "  -3.48%  "
Set-TextValue "D42" "0.0313"
$ws.Range("E42").Value = "  -6.08%  "
Set-TextValue "D43" "3.83"
$ws.Range("E43").Value = "  -6.60%  "
Set-TextValue "D44" "2.026.92"
$ws.Range("E45").Value = "  +0.04%  "
Set-TextValue "D46" "87.19"
$ws.Range("E46").Value = "  -6.85%  "
Set-TextValue "D48" "9.02"
$ws.Range("E48").Value = "  -5.26%  "
Set-TextValue "D49" "2.793.96"
$ws.Range("E49").Value = "  -2.59%  "
Set-TextValue "D50" "104.10"
$ws.Range("E50").Value = "  -4.70%  "
Set-TextValue "D51" "0.190"
$ws.Range("E51").Value = "  -5.52%  "
